# "Generate Report for Handback" -------------------------------------------
# For each localized sheet (zh-cn, de-de) populate the "Latest Target File"
# (F) and "Latest Handback File" (G) columns for the two data rows, stamp
# the "Latest Handback DateTime" (H) column with the handback timestamp, and
# flip every "Status" cell (Overview + both language sheets) from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$HyperlinkFontColor = 15570276   # RGB(0x64,0x95,0xED) == style "HyperLink" (FF6495ED) already used in the workbook
$NewStatusText      = "Handed back: in sync with en-US"

function Set-HandbackLink {
    param($ws, [string]$cellRef, [string]$text, [string]$target)

    $rng = $ws.Range($cellRef)
    $rng.Value = $text
    $ws.Hyperlinks.Add($rng, $target, "", "", $text) | Out-Null
    # Match the existing custom "HyperLink" look (underlined, #6495ED) that
    # the other hyperlink cells on this sheet already use.
    $rng.Font.Underline = $true
    $rng.Font.Color = $HyperlinkFontColor
}

# ---------------------------------------------------------------------------
# Overview sheet: the Status columns (zh-cn / de-de) mirror the same text as
# each language sheet's Status column, so they move to the handback wording
# too.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $NewStatusText
$overview.Range("C2").Value = $NewStatusText
$overview.Range("B3").Value = $NewStatusText
$overview.Range("C3").Value = $NewStatusText

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $NewStatusText
$zhcn.Range("C3").Value = $NewStatusText

Set-HandbackLink $zhcn "F2" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4df332fe9ea659b4cb92ab38a8b30ee284085f02/e2e/a.md"
Set-HandbackLink $zhcn "G2" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/55ea9fd245f451ab3dd10fde5ec13186eb57a0b5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("H2").Value = "2016-03-21 10:29:16"

Set-HandbackLink $zhcn "F3" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4df332fe9ea659b4cb92ab38a8b30ee284085f02/e2e/a.md"
Set-HandbackLink $zhcn "G3" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/55ea9fd245f451ab3dd10fde5ec13186eb57a0b5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("H3").Value = "2016-03-21 10:29:16"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $NewStatusText
$dede.Range("C3").Value = $NewStatusText

Set-HandbackLink $dede "F2" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4df332fe9ea659b4cb92ab38a8b30ee284085f02/e2e/a.md"
Set-HandbackLink $dede "G2" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d626a89116144e1700551e96d4ce5a6d40edf44a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("H2").Value = "2016-03-21 10:29:23"

Set-HandbackLink $dede "F3" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4df332fe9ea659b4cb92ab38a8b30ee284085f02/e2e/a.md"
Set-HandbackLink $dede "G3" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d626a89116144e1700551e96d4ce5a6d40edf44a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("H3").Value = "2016-03-21 10:29:23"
